$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

Set-TextValue "D2" '309.05'
Set-TextValue "E2" '0.40%'
Set-TextValue "D3" '40.92'
Set-TextValue "E3" '-0.13%'
Set-TextValue "D4" '5.126'
Set-TextValue "E4" '1.62%'
Set-TextValue "D5" '0.07634'
Set-TextValue "E5" '-0.06%'
Set-TextValue "D6" '1.607'
Set-TextValue "E6" '-0.33%'
Set-TextValue "D7" '2.472'
Set-TextValue "E7" '1.14%'
Set-TextValue "D8" '0.9087'
Set-TextValue "E8" '0.15%'
Set-TextValue "D9" '0.1273'
Set-TextValue "E9" '26.02%'
Set-TextValue "D10" '0.1803'
Set-TextValue "E10" '1.98%'
Set-TextValue "D11" '0.09059'
Set-TextValue "E11" '-1.45%'
Set-TextValue "D12" '0.04346'
Set-TextValue "E12" '2.86%'
Set-TextValue "E13" '-0.64%'
Set-TextValue "D14" '0.001251'
Set-TextValue "E14" '-1.04%'
Set-TextValue "D15" '0.005654'
Set-TextValue "E15" '-2.96%'
Set-TextValue "D16" '3.354'
Set-TextValue "E16" '-0.05%'
Set-TextValue "D17" '4.293'
Set-TextValue "E17" '0.74%'
Set-TextValue "E18" '1.38%'
Set-TextValue "D19" '6.904'
Set-TextValue "E19" '2.01%'
Set-TextValue "E20" '2.05%'
Set-TextValue "D21" '0.2739'
Set-TextValue "E21" '0.61%'
Set-TextValue "D22" '0.04045'
Set-TextValue "E22" '-2.65%'
Set-TextValue "E23" '4.34%'
Set-TextValue "D24" '0.004051'
Set-TextValue "E24" '-0.60%'
Set-TextValue "E26" '24.75%'
Set-TextValue "D38" '0.02418'
Set-TextValue "E38" '0.26%'
Set-TextValue "D39" '0.05220'
Set-TextValue "E39" '0.95%'
Set-TextValue "D40" '0.007845'
Set-TextValue "E40" '1.21%'
Set-TextValue "D41" '0.1302'
Set-TextValue "E41" '-0.44%'
Set-TextValue "E42" '-4.13%'
Set-TextValue "D43" '0.001842'
Set-TextValue "E43" '-5.42%'
Set-TextValue "D44" '0.007425'
Set-TextValue "D45" '0.3351'
Set-TextValue "E45" '9.67%'
Set-TextValue "D46" '0.00006877'
Set-TextValue "E46" '7.91%'
Set-TextValue "D47" '0.00000000751'
Set-TextValue "E47" '0.18%'
Set-TextValue "E48" '2,455.38%'
Set-TextValue "E49" '-31.75%'
Set-TextValue "D50" '0.00002103'
Set-TextValue "E50" '0.18%'
Set-TextValue "D51" '0.0002003'
Set-TextValue "E51" '0.18%'
